# Template-DFC.xlsx update: new "CEF"/"extcef" and "Informakon"/"rec" data
# points on the "Por empresa" sheet, plus a selection move.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Por empresa")

# Row 3 ("Repasse"): add "extcef" next to the existing "CEF" in C3.
$ws.Range("D3").Value = "extcef"

# Row 4 ("Pró soluto"): add "rec" next to the existing "Informakon" in C4.
$ws.Range("D4").Value = "rec"

# Row 25 ("Captação de empréstimos / financiamentos bancários"): populate
# the previously empty C25/D25 pair with "CEF" / "extcef", matching row 3.
$ws.Range("C25").Value = "CEF"
$ws.Range("D25").Value = "extcef"

# Move the active selection from D16 to B5.
$ws.Range("B5").Select() | Out-Null
